$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Risk_Level" column (D) entirely, shifting cells left
$ws.Range("D1:D3").Delete()

# Add two new risk rows
$ws.Range("A4").Value = "Risk 3"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 3

$ws.Range("A5").Value = "Risk 4"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 1

# Copy the formatting of row 3 (A3:C3) down to the new rows
$ws.Range("A3:C3").Copy()
$ws.Range("A4:C5").PasteSpecial(-4122)

$ws.Range("C6").Select() | Out-Null
